$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1312.4
$ws.Range("I6").Value = 187.33333
$ws.Range("K6").Value = 561.99999
$ws.Range("M6").Value = -449.99999
$ws.Range("H113").Value = 14900.5
$ws.Range("I113").Value = 4002.5
$ws.Range("J113").Value = 17625
$ws.Range("K113").Value = 4002.5
$ws.Range("L113").Value = 17625
$ws.Range("M113").Value = -748.5
$ws.Range("N113").Value = -24133
$ws.Range("H115").Value = 1760.9
$ws.Range("I115").Value = 1776.25
$ws.Range("J115").Value = 1699.5
$ws.Range("K115").Value = 5328.75
$ws.Range("L115").Value = 5098.5
$ws.Range("M115").Value = -3761.75
$ws.Range("N115").Value = -8232.5
$ws.Range("H116").Value = 917571.25
$ws.Range("I116").Value = 2502819.8
$ws.Range("J116").Value = 11715
$ws.Range("K116").Value = 2502819.8
$ws.Range("L116").Value = 11715
$ws.Range("M116").Value = -2499377.8
$ws.Range("N116").Value = -18599
$ws.Range("H132").Value = 35721052
$ws.Range("I132").Value = 40006536
$ws.Range("J132").Value = 8666.666999999999
$ws.Range("K132").Value = 120019608
$ws.Range("L132").Value = 26000.001
$ws.Range("M132").Value = -120017078
$ws.Range("N132").Value = -31060.001
$ws.Range("H137").Value = 3615.7
$ws.Range("I137").Value = 3358.7878
$ws.Range("J137").Value = 4114.4116
$ws.Range("K137").Value = 10076.3634
$ws.Range("L137").Value = 12343.2348
$ws.Range("M137").Value = -7526.3634
$ws.Range("N137").Value = -17443.2348
$ws.Range("H138").Value = 3631.4177
$ws.Range("I138").Value = 1469.7142
$ws.Range("J138").Value = 4414.1035
$ws.Range("K138").Value = 4409.142599999999
$ws.Range("L138").Value = 13242.3105
$ws.Range("M138").Value = 730.8574000000008
$ws.Range("N138").Value = -23522.3105
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11997.543
$ws.Range("I32").Value = 8401.940000000001
$ws.Range("J32").Value = 17796.902
$ws.Range("K32").Value = 8401.940000000001
$ws.Range("L32").Value = 17796.902
$ws.Range("M32").Value = -8114.940000000001
$ws.Range("N32").Value = -18370.902
$ws.Range("H45").Value = 1077.6923
$ws.Range("I45").Value = 972.8570999999999
$ws.Range("K45").Value = 972.8570999999999
$ws.Range("M45").Value = -595.8570999999999
$ws.Range("H110").Value = 561.86365
$ws.Range("I110").Value = 536.7222
$ws.Range("J110").Value = 675
$ws.Range("K110").Value = 536.7222
$ws.Range("L110").Value = 675
$ws.Range("M110").Value = 1508.2778
$ws.Range("N110").Value = -4765
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1906.12
$ws.Range("I105").Value = 1831.2354
$ws.Range("J105").Value = 2065.25
$ws.Range("K105").Value = 1831.2354
$ws.Range("L105").Value = 2065.25
$ws.Range("M105").Value = -84.23540000000003
$ws.Range("N105").Value = -5559.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1371.9584
$ws.Range("I16").Value = 1159.2
$ws.Range("K16").Value = 1159.2
$ws.Range("M16").Value = -872.2
$ws.Range("H31").Value = 4867.6523
$ws.Range("I31").Value = 2410.2222
$ws.Range("J31").Value = 6447.4287
$ws.Range("K31").Value = 2410.2222
$ws.Range("L31").Value = 6447.4287
$ws.Range("M31").Value = -2115.2222
$ws.Range("N31").Value = -7037.4287
$ws.Range("H34").Value = 4867.6523
$ws.Range("I34").Value = 2410.2222
$ws.Range("J34").Value = 6447.4287
$ws.Range("K34").Value = 2410.2222
$ws.Range("L34").Value = 6447.4287
$ws.Range("M34").Value = -2208.2222
$ws.Range("N34").Value = -6851.4287
$ws.Range("H58").Value = 2348.1553
$ws.Range("I58").Value = 1822.5283
$ws.Range("J58").Value = 7919.8
$ws.Range("K58").Value = 1822.5283
$ws.Range("L58").Value = 7919.8
$ws.Range("M58").Value = -1619.5283
$ws.Range("N58").Value = -8325.799999999999
$ws.Range("H113").Value = 1371.9584
$ws.Range("I113").Value = 1159.2
$ws.Range("K113").Value = 1159.2
$ws.Range("M113").Value = 1010.8
$ws.Range("H136").Value = 2348.1553
$ws.Range("I136").Value = 1822.5283
$ws.Range("J136").Value = 7919.8
$ws.Range("K136").Value = 5467.5849
$ws.Range("L136").Value = 23759.4
$ws.Range("M136").Value = -2917.5849
$ws.Range("N136").Value = -28859.4
$ws.Range("H140").Value = 143543.33
$ws.Range("J140").Value = 143543.33
$ws.Range("L140").Value = 143543.33
$ws.Range("N140").Value = -153903.33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 8825.058999999999
$ws.Range("I26").Value = 17299.715
$ws.Range("K26").Value = 51899.145
$ws.Range("M26").Value = -51611.145
$ws.Range("H86").Value = 15133.833
$ws.Range("I86").Value = 16160
$ws.Range("K86").Value = 48480
$ws.Range("M86").Value = -47294
$ws.Range("H89").Value = 15133.833
$ws.Range("I89").Value = 16160
$ws.Range("K89").Value = 145440
$ws.Range("M89").Value = -139512
$ws.Range("H94").Value = 2458.8235
$ws.Range("J94").Value = 2720
$ws.Range("L94").Value = 8160
$ws.Range("N94").Value = -9512
$ws.Range("H97").Value = 538.5417
$ws.Range("J97").Value = 582
$ws.Range("L97").Value = 1746
$ws.Range("N97").Value = -2738
$ws.Range("H113").Value = 541.1702
$ws.Range("I113").Value = 555.5238000000001
$ws.Range("J113").Value = 529.5769
$ws.Range("K113").Value = 1666.5714
$ws.Range("L113").Value = 1588.7307
$ws.Range("M113").Value = 503.4285999999997
$ws.Range("N113").Value = -5928.7307
$ws.Range("H118").Value = 8641.286
$ws.Range("I118").Value = 489
$ws.Range("J118").Value = 10000
$ws.Range("K118").Value = 1467
$ws.Range("L118").Value = 30000
$ws.Range("M118").Value = -224
$ws.Range("N118").Value = -32486
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33333.332
$ws.Range("J52").Value = 33333.332
$ws.Range("L52").Value = 33333.332
$ws.Range("N52").Value = -33851.332
$ws.Range("H97").Value = 2902.2222
$ws.Range("I97").Value = 2024
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 2024
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -1528
$ws.Range("N97").Value = -4992
$ws.Range("H107").Value = 450.3793
$ws.Range("I107").Value = 198.73685
$ws.Range("J107").Value = 928.5
$ws.Range("K107").Value = 198.73685
$ws.Range("L107").Value = 928.5
$ws.Range("M107").Value = 1721.26315
$ws.Range("N107").Value = -4768.5
$ws.Range("H132").Value = 2600.7222
$ws.Range("I132").Value = 580.6875
$ws.Range("J132").Value = 4216.75
$ws.Range("K132").Value = 1742.0625
$ws.Range("L132").Value = 12650.25
$ws.Range("M132").Value = 787.9375
$ws.Range("N132").Value = -17710.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13891478
$ws.Range("I22").Value = 25002070
$ws.Range("J22").Value = 3237.875
$ws.Range("K22").Value = 25002070
$ws.Range("L22").Value = 3237.875
$ws.Range("M22").Value = -25001775
$ws.Range("N22").Value = -3827.875
$ws.Range("H27").Value = 13891478
$ws.Range("I27").Value = 25002070
$ws.Range("J27").Value = 3237.875
$ws.Range("K27").Value = 25002070
$ws.Range("L27").Value = 3237.875
$ws.Range("M27").Value = -25001963
$ws.Range("N27").Value = -3451.875
$ws.Range("H93").Value = 2133.7727
$ws.Range("I93").Value = 1865.421
$ws.Range("K93").Value = 1865.421
$ws.Range("M93").Value = -617.421
$ws.Range("H127").Value = 32055.416
$ws.Range("J127").Value = 32055.416
$ws.Range("L127").Value = 32055.416
$ws.Range("N127").Value = -41975.416
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 33353200
$ws.Range("J52").Value = 29799.5
$ws.Range("L52").Value = 29799.5
$ws.Range("N52").Value = -30251.5
$ws.Range("H81").Value = 40179844
$ws.Range("I81").Value = 40179844
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 80359688
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -80358627
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 40179844
$ws.Range("I84").Value = 40179844
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 401798440
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -401793136
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 12826597
$ws.Range("I132").Value = 7999.7144
$ws.Range("J132").Value = 27781626
$ws.Range("K132").Value = 23999.1432
$ws.Range("L132").Value = 83344878
$ws.Range("M132").Value = -21469.1432
$ws.Range("N132").Value = -83349938
$ws.Range("H136").Value = 3478.5925
$ws.Range("I136").Value = 800.2353000000001
$ws.Range("K136").Value = 2400.7059
$ws.Range("M136").Value = 149.2941000000001
